$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3815.4639
$ws.Range("I15").Value = 3815.4639
$ws.Range("K15").Value = 11446.3917
$ws.Range("M15").Value = -11277.3917

$ws.Range("H19").Value = 793.26086
$ws.Range("I19").Value = 322.81818
$ws.Range("J19").Value = 1224.5
$ws.Range("K19").Value = 322.81818
$ws.Range("L19").Value = 1224.5
$ws.Range("M19").Value = -147.81818
$ws.Range("N19").Value = -1574.5

$ws.Range("H53").Value = 1719.6428
$ws.Range("I53").Value = 1358.5
$ws.Range("J53").Value = 2622.5
$ws.Range("K53").Value = 1358.5
$ws.Range("L53").Value = 2622.5
$ws.Range("M53").Value = -721.5
$ws.Range("N53").Value = -3896.5

$ws.Range("H103").Value = 1423.3334
$ws.Range("I103").Value = 770
$ws.Range("K103").Value = 2310
$ws.Range("M103").Value = -1724

$ws.Range("H116").Value = 3113.95
$ws.Range("I116").Value = 2421.1538
$ws.Range("J116").Value = 4400.5713
$ws.Range("K116").Value = 2421.1538
$ws.Range("L116").Value = 4400.5713
$ws.Range("M116").Value = 1020.8462
$ws.Range("N116").Value = -11284.5713

$ws.Range("H129").Value = 860.4706
$ws.Range("J129").Value = 881.9583
$ws.Range("L129").Value = 2645.8749
$ws.Range("N129").Value = -12645.8749

$ws.Range("H137").Value = 2359.5
$ws.Range("I137").Value = 1498
$ws.Range("K137").Value = 4494
$ws.Range("M137").Value = -1944

$ws.Range("H138").Value = 1951.04
$ws.Range("I138").Value = 774.4666999999999
$ws.Range("J138").Value = 2158.6707
$ws.Range("K138").Value = 2323.4001
$ws.Range("L138").Value = 6476.0121
$ws.Range("M138").Value = 2816.5999
$ws.Range("N138").Value = -16756.0121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4069.9302
$ws.Range("I32").Value = 4000.6829
$ws.Range("K32").Value = 4000.6829
$ws.Range("M32").Value = -3713.6829

$ws.Range("H122").Value = 1317.909
$ws.Range("I122").Value = 1266.4445
$ws.Range("K122").Value = 3799.3335
$ws.Range("M122").Value = -1349.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null

$ws.Range("H88").Value = 31332.166
$ws.Range("J88").Value = 31332.166
$ws.Range("L88").Value = 31332.166
$ws.Range("N88").Value = -32144.166

$ws.Range("H91").Value = 31332.166
$ws.Range("J91").Value = 31332.166
$ws.Range("L91").Value = 31332.166
$ws.Range("N91").Value = -34140.166

$ws.Range("H123").Value = 36567.2
$ws.Range("I123").Value = 35709
$ws.Range("J123").Value = 40000
$ws.Range("K123").Value = 35709
$ws.Range("L123").Value = 40000
$ws.Range("M123").Value = -30809
$ws.Range("N123").Value = -49800

$ws.Range("H132").Value = 20571
$ws.Range("J132").Value = 20571
$ws.Range("L132").Value = 20571
$ws.Range("N132").Value = -30691

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

$ws.Range("H140").Value = 22853.912
$ws.Range("J140").Value = 22948.182
$ws.Range("L140").Value = 22948.182
$ws.Range("N140").Value = -33308.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1075.8572
$ws.Range("I58").Value = 1084.3182
$ws.Range("J58").Value = 1061.5385
$ws.Range("K58").Value = 1084.3182
$ws.Range("L58").Value = 1061.5385
$ws.Range("M58").Value = -881.3181999999999
$ws.Range("N58").Value = -1467.5385

$ws.Range("H122").Value = 621
$ws.Range("I122").Value = 478
$ws.Range("J122").Value = 907
$ws.Range("K122").Value = 1434
$ws.Range("L122").Value = 2721
$ws.Range("M122").Value = 1016
$ws.Range("N122").Value = -7621

$ws.Range("H136").Value = 1075.8572
$ws.Range("I136").Value = 1084.3182
$ws.Range("J136").Value = 1061.5385
$ws.Range("K136").Value = 3252.9546
$ws.Range("L136").Value = 3184.6155
$ws.Range("M136").Value = -702.9546
$ws.Range("N136").Value = -8284.6155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2364.875
$ws.Range("I68").Value = 841.3333
$ws.Range("K68").Value = 2523.9999
$ws.Range("M68").Value = -1712.9999

$ws.Range("H71").Value = 2364.875
$ws.Range("I71").Value = 841.3333
$ws.Range("K71").Value = 7571.9997
$ws.Range("M71").Value = -3515.9997

$ws.Range("H122").Value = 1100.5
$ws.Range("I122").Value = 699.6667
$ws.Range("J122").Value = 1501.3334
$ws.Range("K122").Value = 6297.0003
$ws.Range("L122").Value = 13512.0006
$ws.Range("M122").Value = -3847.0003
$ws.Range("N122").Value = -18412.0006

$ws.Range("H131").Value = 20001406
$ws.Range("J131").Value = 1512.159
$ws.Range("L131").Value = 4536.477000000001
$ws.Range("N131").Value = -14616.477

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").Value = $null

$ws.Range("H113").Value = 1576.6666
$ws.Range("I113").Value = 1576.6666
$ws.Range("K113").Value = 1576.6666
$ws.Range("M113").Value = 593.3334

$ws.Range("H122").Value = 2876.8462
$ws.Range("I122").Value = 1933.2222
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 5799.6666
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3349.6666
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2562.8572
$ws.Range("I40").Value = 2322.9167
$ws.Range("K40").Value = 2322.9167
$ws.Range("M40").Value = -2186.9167

$ws.Range("H100").Value = 2501.6667
$ws.Range("I100").Value = 2500
$ws.Range("J100").Value = 2502.5
$ws.Range("K100").Value = 2500
$ws.Range("L100").Value = 2502.5
$ws.Range("M100").Value = -1959
$ws.Range("N100").Value = -3584.5

$ws.Range("H122").Value = 10122788
$ws.Range("J122").Value = 3995.4614
$ws.Range("L122").Value = 11986.3842
$ws.Range("N122").Value = -16886.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3666.6667

$ws.Range("H122").Value = 14446580
$ws.Range("I122").Value = 16252190
$ws.Range("J122").Value = 1702.5
$ws.Range("K122").Value = 48756570
$ws.Range("L122").Value = 5107.5
$ws.Range("M122").Value = -48754120
$ws.Range("N122").Value = -10007.5

$ws.Range("H139").Value = 41926
$ws.Range("J139").Value = 46543.332
$ws.Range("L139").Value = 46543.332
$ws.Range("N139").Value = -56823.332

$ws.Range("H141").Value = 44085.832
$ws.Range("J141").Value = 44085.832
$ws.Range("L141").Value = 44085.832
$ws.Range("N141").Value = -54445.832
